# fix(publipostage): Correct status name
#
# "bleu" -> "noir"
# "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
# "résultat et / ou publication posté" -> "résultat postés ou publiés"
# "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
# "résultat et / ou publication posté dans les 12 mois" -> "résultat postés ou publiés dans les 12 mois"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$whole = [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole

$ws.Cells.Replace("bleu", "noir", $whole)

# Replace the longest / most specific strings first so that the shorter
# "résultat et / ou publication posté" text isn't matched as a prefix
# (xlWhole protects against this too, but keep the safer ordering).
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", $whole)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", $whole)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", $whole)
